$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.274.76"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "'1.576.01"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'207.93"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  -1.87%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'1.800.93"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Value = "'1.577.81"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "'3.78"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").Value = "'27.297.98"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "'62.47"
$ws.Range("E17").Value = "  -1.12%  "
$ws.Range("D18").Value = "'215.10"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "'7.34"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "'0.0₃0688"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'9.41"
$ws.Range("E23").Value = "  -3.53%  "
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("D25").Value = "'151.93"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  -3.76%  "
$ws.Range("D27").Value = "'14.97"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.104"
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("E30").Value = "  -1.61%  "
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").Value = "'1.410.38"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").Value = "'0.939"
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").Value = "'0.519"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("E42").Value = "  +2.28%  "
$ws.Range("D43").Value = "'1.82"
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("D45").Value = "'63.95"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D47").Value = "'1.713.37"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "'86.13"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "'0.0₇0985"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "'0.0953"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").Value = "'0.0494"
$ws.Range("E51").Value = "  -0.06%  "
